$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A few "Price" values are purely numeric-looking strings (e.g. "1.620")
# that Excel would otherwise auto-convert to a Number and normalize
# (dropping the trailing zero). Force those specific cells to Text first
# so the literal string is preserved, then restore default formatting.
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

$ws.Range("D2").Value = '27.160.96'
$ws.Range("E2").Value = '  -1.19%  '
$ws.Range("D3").Value = '1.782.16'
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").Value = '336.47'
$ws.Range("E5").Value = '  -2.06%  '
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("D7").Value = '0.3889'
$ws.Range("E7").Value = '  +1.14%  '
$ws.Range("E8").Value = '  -2.90%  '
$ws.Range("D9").Value = '47.81'
$ws.Range("E9").Value = '  -2.58%  '
$ws.Range("D10").Value = '1.185'
$ws.Range("E10").Value = '  -4.16%  '
$ws.Range("D11").Value = '0.07419'
$ws.Range("D12").Value = '1.003'
$ws.Range("E12").Value = '  +0.14%  '
$ws.Range("D13").Value = '21.53'
$ws.Range("E13").Value = '  -3.81%  '
$ws.Range("D14").Value = '6.421'
$ws.Range("E14").Value = '  -2.84%  '
$ws.Range("D15").Value = '1.780.57'
$ws.Range("E15").Value = '  -2.22%  '
$ws.Range("D16").Value = '7.083'
$ws.Range("E16").Value = '  -1.96%  '
$ws.Range("D17").Value = '0.00001086'
$ws.Range("E17").Value = '  -3.23%  '
$ws.Range("E18").Value = '  -1.03%  '
$ws.Range("D19").Value = '83.17'
$ws.Range("E19").Value = '  -3.81%  '
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.11%  '
$ws.Range("D21").Value = '17.54'
$ws.Range("E21").Value = '  -0.80%  '
$ws.Range("D22").Value = '6.478'
$ws.Range("E22").Value = '  -1.49%  '
$ws.Range("D23").Value = '27.163.55'
$ws.Range("E23").Value = '  -1.22%  '
$ws.Range("E24").Value = '  -6.90%  '
$ws.Range("D25").Value = '2.374'
$ws.Range("E25").Value = '  -3.74%  '
$ws.Range("D26").Value = '21.06'
$ws.Range("E26").Value = '  -5.83%  '
$ws.Range("D27").Value = '2.481'
$ws.Range("E27").Value = '  -7.91%  '
$ws.Range("D28").Value = '1.436'
$ws.Range("E28").Value = '  -2.50%  '
$ws.Range("D29").Value = '155.81'
$ws.Range("E29").Value = '  +1.16%  '
$ws.Range("D30").Value = '1.982.59'
$ws.Range("E30").Value = '  -2.07%  '
$ws.Range("D31").Value = '133.89'
$ws.Range("E31").Value = '  -1.91%  '
$ws.Range("D32").Value = '3.973'
$ws.Range("E32").Value = '  -2.26%  '
$ws.Range("D33").Value = '5.939'
$ws.Range("E33").Value = '  -7.07%  '
$ws.Range("D34").Value = '0.08706'
$ws.Range("E34").Value = '  -1.27%  '
$ws.Range("D35").Value = '12.82'
$ws.Range("E35").Value = '  -8.14%  '
$ws.Range("D36").Value = '1.620'
$ws.Range("E36").Value = '  -3.99%  '
$ws.Range("D37").Value = '5.385'
$ws.Range("E37").Value = '  -4.33%  '
$ws.Range("D38").Value = '0.6763'
$ws.Range("E38").Value = '  -4.16%  '
$ws.Range("D39").Value = '0.06327'
$ws.Range("E39").Value = '  -2.47%  '
$ws.Range("E40").Value = '  -2.97%  '
$ws.Range("D41").Value = '0.2181'
$ws.Range("E41").Value = '  -3.85%  '
$ws.Range("D42").Value = '1.239'
$ws.Range("E42").Value = '  -4.51%  '
$ws.Range("D43").Value = '8.401'
$ws.Range("E43").Value = '  -6.35%  '
$ws.Range("D44").Value = '14.17'
$ws.Range("E44").Value = '  -4.62%  '
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("D46").Value = '0.6357'
$ws.Range("E46").Value = '  -3.87%  '
$ws.Range("D47").Value = '3.845'
$ws.Range("E47").Value = '  -2.79%  '
$ws.Range("D48").Value = '2.140'
$ws.Range("E48").Value = '  -2.33%  '
$ws.Range("D49").Value = '130.59'
$ws.Range("E49").Value = '  -1.66%  '
$ws.Range("D50").Value = '0.07116'
$ws.Range("E50").Value = '  -3.06%  '
$ws.Range("D51").Value = '78.89'
$ws.Range("E51").Value = '  -2.32%  '

# Restore the cells we forced to Text back to the default (unstyled) format
# so only the cell VALUES differ from the original, not their formatting.
$ws.Range("D36").ClearFormats()
$ws.Range("D48").ClearFormats()
